$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Rows 2-12, column A ("card") change from "2" to "24".
# Force the values to be written as text (matching the existing
# inline-string cell type) rather than being auto-coerced to a number,
# then strip the quote-prefix formatting it introduces so the cell
# keeps its original (default) style.
for ($r = 2; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.NumberFormat = "@"
    $cell.Value = "24"
}
$ws.Range("A2:A12").ClearFormats()
